# Add a space before ":" in the "statut_name" column (column B) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 1; $r -le $rows; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value()
    if ($val -ne $null -and $val -is [string] -and $val -match '^(\d):\s') {
        $newVal = $val -replace '^(\d):\s', '$1 : '
        $cell.Value = $newVal
    }
}
